$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 10 new match results (Giornata 8) as rows 72-81, following the existing
# data which ends at row 71 (dimension grows from A1:H71 to A1:H81).

$ws.Cells.Item(72, 1).Value = 70
$ws.Cells.Item(72, 2).Value = "CROTONE"
$ws.Cells.Item(72, 3).Value = "LAZIO"
$ws.Cells.Item(72, 4).Value = -1
$ws.Cells.Item(72, 5).Value = 4.75
$ws.Cells.Item(72, 6).Value = 4
$ws.Cells.Item(72, 7).Value = 1.66
$ws.Cells.Item(72, 8).Value = 8

$ws.Cells.Item(73, 1).Value = 71
$ws.Cells.Item(73, 2).Value = "SPEZIA"
$ws.Cells.Item(73, 3).Value = "ATALANTA"
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 6.5
$ws.Cells.Item(73, 6).Value = 5.25
$ws.Cells.Item(73, 7).Value = 1.4
$ws.Cells.Item(73, 8).Value = 8

$ws.Cells.Item(74, 1).Value = 72
$ws.Cells.Item(74, 2).Value = "JUVENTUS"
$ws.Cells.Item(74, 3).Value = "CAGLIARI"
$ws.Cells.Item(74, 4).Value = 1
$ws.Cells.Item(74, 5).Value = 1.28
$ws.Cells.Item(74, 6).Value = 5.5
$ws.Cells.Item(74, 7).Value = 10
$ws.Cells.Item(74, 8).Value = 8

$ws.Cells.Item(75, 1).Value = 73
$ws.Cells.Item(75, 2).Value = "FIORENTINA"
$ws.Cells.Item(75, 3).Value = "BENEVENTO"
$ws.Cells.Item(75, 4).Value = -1
$ws.Cells.Item(75, 5).Value = 1.5
$ws.Cells.Item(75, 6).Value = 4.75
$ws.Cells.Item(75, 7).Value = 5.5
$ws.Cells.Item(75, 8).Value = 8

$ws.Cells.Item(76, 1).Value = 74
$ws.Cells.Item(76, 2).Value = "INTER"
$ws.Cells.Item(76, 3).Value = "TORINO"
$ws.Cells.Item(76, 4).Value = 1
$ws.Cells.Item(76, 5).Value = 1.33
$ws.Cells.Item(76, 6).Value = 5.75
$ws.Cells.Item(76, 7).Value = 8
$ws.Cells.Item(76, 8).Value = 8

$ws.Cells.Item(77, 1).Value = 75
$ws.Cells.Item(77, 2).Value = "ROMA"
$ws.Cells.Item(77, 3).Value = "PARMA"
$ws.Cells.Item(77, 4).Value = 1
$ws.Cells.Item(77, 5).Value = 1.5
$ws.Cells.Item(77, 6).Value = 4.75
$ws.Cells.Item(77, 7).Value = 5.75
$ws.Cells.Item(77, 8).Value = 8

$ws.Cells.Item(78, 1).Value = 76
$ws.Cells.Item(78, 2).Value = "SAMPDORIA"
$ws.Cells.Item(78, 3).Value = "BOLOGNA"
$ws.Cells.Item(78, 4).Value = -1
$ws.Cells.Item(78, 5).Value = 2.5
$ws.Cells.Item(78, 6).Value = 3.5
$ws.Cells.Item(78, 7).Value = 2.7
$ws.Cells.Item(78, 8).Value = 8

$ws.Cells.Item(79, 1).Value = 77
$ws.Cells.Item(79, 2).Value = "VERONA"
$ws.Cells.Item(79, 3).Value = "SASSUOLO"
$ws.Cells.Item(79, 4).Value = -1
$ws.Cells.Item(79, 5).Value = 3.2
$ws.Cells.Item(79, 6).Value = 3.8
$ws.Cells.Item(79, 7).Value = 2.1
$ws.Cells.Item(79, 8).Value = 8

$ws.Cells.Item(80, 1).Value = 78
$ws.Cells.Item(80, 2).Value = "UDINESE"
$ws.Cells.Item(80, 3).Value = "GENOA"
$ws.Cells.Item(80, 4).Value = 1
$ws.Cells.Item(80, 5).Value = 1.8
$ws.Cells.Item(80, 6).Value = 3.75
$ws.Cells.Item(80, 7).Value = 4.33
$ws.Cells.Item(80, 8).Value = 8

$ws.Cells.Item(81, 1).Value = 79
$ws.Cells.Item(81, 2).Value = "NAPOLI"
$ws.Cells.Item(81, 3).Value = "MILAN"
$ws.Cells.Item(81, 4).Value = -1
$ws.Cells.Item(81, 5).Value = 2.3
$ws.Cells.Item(81, 6).Value = 3.5
$ws.Cells.Item(81, 7).Value = 3
$ws.Cells.Item(81, 8).Value = 8

# Match the bold/centered/bordered style used by column A in the existing data rows
# (A2:A71) by copying its formatting onto the new A72:A81 cells.
$ws.Range("A2").Copy()
$ws.Range("A72:A81").PasteSpecial(-4122)
$excel.CutCopyMode = $false

